$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 47

# Column A must hold the literal text "2020-07-16", matching the other
# date cells in the column (they are stored as shared-string text, not
# real Excel dates). Assigning the string straight to the cell makes
# Excel's input auto-recognition convert it into a date serial number, so
# instead stage it as a formula in a scratch cell far outside the used
# range (a formula result is never re-parsed as a literal date), copy just
# the computed value over to the target cell, then clear the scratch cell
# completely so nothing else is left behind.
$scratch = $ws.Cells.Item(1000, 100)
$scratch.Formula = '="2020-07-16"'
$scratch.Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163) | Out-Null
$scratch.Clear()

$ws.Cells.Item($row, 2).Value = 324041
$ws.Cells.Item($row, 3).Value = 375455
$ws.Cells.Item($row, 4).Value = 82567
$ws.Cells.Item($row, 5).Value = 37574
$ws.Cells.Item($row, 6).Value = 28.95
